$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Fgf1"
$ws.Cells.Item(2, 3).Value = "Fgfr4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 2.004760666666666
$ws.Cells.Item(2, 8).Value = 6.014282
$ws.Cells.Item(2, 9).Value = 0.1200698528618338
$ws.Cells.Item(2, 10).Value = 0.1200698528618338
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.5555316666666666
$ws.Cells.Item(2, 14).Value = 1.666595
$ws.Cells.Item(2, 15).Value = 0.01938483203642842
$ws.Cells.Item(2, 16).Value = 0.01938483203642843
$ws.Cells.Item(2, 17).Value = 1.113708034421111
$ws.Cells.Item(2, 18).Value = 10.02337230979
$ws.Cells.Item(2, 19).Value = 0.002327533930365323
$ws.Cells.Item(2, 20).Value = 0.002327533930365323

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Fgf1"
$ws.Cells.Item(3, 3).Value = "Fgfr4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 2.004760666666666
$ws.Cells.Item(3, 8).Value = 6.014282
$ws.Cells.Item(3, 9).Value = 0.1200698528618338
$ws.Cells.Item(3, 10).Value = 0.1200698528618338
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.1817723333333333
$ws.Cells.Item(3, 14).Value = 0.545317
$ws.Cells.Item(3, 15).Value = 0.006342799811357313
$ws.Cells.Item(3, 16).Value = 0.006342799811357313
$ws.Cells.Item(3, 17).Value = 0.3644100241548889
$ws.Cells.Item(3, 18).Value = 3.279690217394
$ws.Cells.Item(3, 19).Value = 0.0007615790400817395
$ws.Cells.Item(3, 20).Value = 0.0007615790400817397

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Fgf1"
$ws.Cells.Item(4, 3).Value = "Fgfr4"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 2.004760666666666
$ws.Cells.Item(4, 8).Value = 6.014282
$ws.Cells.Item(4, 9).Value = 0.1200698528618338
$ws.Cells.Item(4, 10).Value = 0.1200698528618338
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 27.92075533333333
$ws.Cells.Item(4, 14).Value = 83.762266
$ws.Cells.Item(4, 15).Value = 0.9742723681522142
$ws.Cells.Item(4, 16).Value = 0.9742723681522143
$ws.Cells.Item(4, 17).Value = 55.97443207589021
$ws.Cells.Item(4, 18).Value = 503.769888683012
$ws.Cells.Item(4, 19).Value = 0.1169807398913867
$ws.Cells.Item(4, 20).Value = 0.1169807398913867

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Fgf1"
$ws.Cells.Item(5, 3).Value = "Fgfr4"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 9.409654999999999
$ws.Cells.Item(5, 8).Value = 28.228965
$ws.Cells.Item(5, 9).Value = 0.5635664696121425
$ws.Cells.Item(5, 10).Value = 0.5635664696121424
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.5555316666666666
$ws.Cells.Item(5, 14).Value = 1.666595
$ws.Cells.Item(5, 15).Value = 0.01938483203642842
$ws.Cells.Item(5, 16).Value = 0.01938483203642843
$ws.Cells.Item(5, 17).Value = 5.227361324908332
$ws.Cells.Item(5, 18).Value = 47.046251924175
$ws.Cells.Item(5, 19).Value = 0.01092464135479433
$ws.Cells.Item(5, 20).Value = 0.01092464135479433

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Fgf1"
$ws.Cells.Item(6, 3).Value = "Fgfr4"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 9.409654999999999
$ws.Cells.Item(6, 8).Value = 28.228965
$ws.Cells.Item(6, 9).Value = 0.5635664696121425
$ws.Cells.Item(6, 10).Value = 0.5635664696121424
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.1817723333333333
$ws.Cells.Item(6, 14).Value = 0.545317
$ws.Cells.Item(6, 15).Value = 0.006342799811357313
$ws.Cells.Item(6, 16).Value = 0.006342799811357313
$ws.Cells.Item(6, 17).Value = 1.710414945211667
$ws.Cells.Item(6, 18).Value = 15.393734506905
$ws.Cells.Item(6, 19).Value = 0.003574589297143204
$ws.Cells.Item(6, 20).Value = 0.003574589297143204

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Fgf1"
$ws.Cells.Item(7, 3).Value = "Fgfr4"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 9.409654999999999
$ws.Cells.Item(7, 8).Value = 28.228965
$ws.Cells.Item(7, 9).Value = 0.5635664696121425
$ws.Cells.Item(7, 10).Value = 0.5635664696121424
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 27.92075533333333
$ws.Cells.Item(7, 14).Value = 83.762266
$ws.Cells.Item(7, 15).Value = 0.9742723681522142
$ws.Cells.Item(7, 16).Value = 0.9742723681522143
$ws.Cells.Item(7, 17).Value = 262.7246750260766
$ws.Cells.Item(7, 18).Value = 2364.52207523469
$ws.Cells.Item(7, 19).Value = 0.5490672389602049
$ws.Cells.Item(7, 20).Value = 0.5490672389602049

$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Fgf1"
$ws.Cells.Item(8, 3).Value = "Fgfr4"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 5.282203999999999
$ws.Cells.Item(8, 8).Value = 15.846612
$ws.Cells.Item(8, 9).Value = 0.3163636775260238
$ws.Cells.Item(8, 10).Value = 0.3163636775260238
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.5555316666666666
$ws.Cells.Item(8, 14).Value = 1.666595
$ws.Cells.Item(8, 15).Value = 0.01938483203642842
$ws.Cells.Item(8, 16).Value = 0.01938483203642843
$ws.Cells.Item(8, 17).Value = 2.934431591793333
$ws.Cells.Item(8, 18).Value = 26.40988432614
$ws.Cells.Item(8, 19).Value = 0.006132656751268777
$ws.Cells.Item(8, 20).Value = 0.006132656751268778

$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Fgf1"
$ws.Cells.Item(9, 3).Value = "Fgfr4"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 5.282203999999999
$ws.Cells.Item(9, 8).Value = 15.846612
$ws.Cells.Item(9, 9).Value = 0.3163636775260238
$ws.Cells.Item(9, 10).Value = 0.3163636775260238
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.1817723333333333
$ws.Cells.Item(9, 14).Value = 0.545317
$ws.Cells.Item(9, 15).Value = 0.006342799811357313
$ws.Cells.Item(9, 16).Value = 0.006342799811357313
$ws.Cells.Item(9, 17).Value = 0.9601585462226666
$ws.Cells.Item(9, 18).Value = 8.641426916004
$ws.Cells.Item(9, 19).Value = 0.002006631474132369
$ws.Cells.Item(9, 20).Value = 0.00200663147413237

$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Fgf1"
$ws.Cells.Item(10, 3).Value = "Fgfr4"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 5.282203999999999
$ws.Cells.Item(10, 8).Value = 15.846612
$ws.Cells.Item(10, 9).Value = 0.3163636775260238
$ws.Cells.Item(10, 10).Value = 0.3163636775260238
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 27.92075533333333
$ws.Cells.Item(10, 14).Value = 83.762266
$ws.Cells.Item(10, 15).Value = 0.9742723681522142
$ws.Cells.Item(10, 16).Value = 0.9742723681522143
$ws.Cells.Item(10, 17).Value = 147.4831255047546
$ws.Cells.Item(10, 18).Value = 1327.348129542792
$ws.Cells.Item(10, 19).Value = 0.3082243893006226
$ws.Cells.Item(10, 20).Value = 0.3082243893006227
